$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '83.414.72'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +4.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.220.59'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +6.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '624.00'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.309'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +26.41%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.589'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.219.39'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.601'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000271'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.11%  '
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.35'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.793.11'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.70'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '83.047.46'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.204.87'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.27'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +9.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.29'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '444.69'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.09'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.22'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.33'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.96%  '
$ws.Range("E25").Value = '  +8.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.85'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +9.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '78.40'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.355.08'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000124'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.23'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.46%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '572.30'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.149'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +21.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.154'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.26%  '
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.03'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.26'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +12.78%  '
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.408'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.82%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.05'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +12.31%  '
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.89'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.03'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +14.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '159.99'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '188.49'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.76'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.33'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.779'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.13'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.19%  '
